$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Invalid (G) and Absent (H) set to 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Total Attendance Count (D) and Real (E) set to 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Rows 5-18: Absent (H) set to 1
for ($r = 5; $r -le 18; $r++) {
    $ws.Range("H$r").Value = 1
}
